$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Athletico-PR vs Atletico-MG) - Odd_Over25_FT / Odd_Under25_FT
$ws.Range("Q2").Value = 2.25
$ws.Range("R2").Value = 1.62

# Row 4 (Amazonas vs Goias) - Odd_Over05_FT / Odd_Under05_FT
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 10

# Row 6 (Ponte Preta vs Sport Recife)
$ws.Range("G6").Value = 4.5
$ws.Range("H6").Value = 3.25
$ws.Range("I6").Value = 1.9
$ws.Range("J6").Value = 4.75
$ws.Range("K6").Value = 2.05
$ws.Range("L6").Value = 2.6
$ws.Range("Q6").Value = 2.15
$ws.Range("R6").Value = 1.67
$ws.Range("X6").Value = 21
$ws.Range("Z6").Value = 41
$ws.Range("AI6").Value = 8.5
$ws.Range("AO6").Value = 23

# Row 9 (Chrobry Glogow vs Wisla)
$ws.Range("G9").Value = 6.5
$ws.Range("H9").Value = 4.75
$ws.Range("I9").Value = 1.36
$ws.Range("J9").Value = 7
$ws.Range("L9").Value = 1.83
$ws.Range("M9").Value = 1.03
$ws.Range("N9").Value = 17
$ws.Range("U9").Value = 1.8
$ws.Range("V9").Value = 1.91
$ws.Range("AE9").Value = 19
$ws.Range("AG9").Value = 201
$ws.Range("AH9").Value = 9
$ws.Range("AK9").Value = 9.5
$ws.Range("AO9").Value = 34
$ws.Range("AP9").Value = 34
$ws.Range("AQ9").Value = 126
$ws.Range("AR9").Value = 126
$ws.Range("AS9").Value = 201
$ws.Range("AU9").Value = 8.5
$ws.Range("AW9").Value = 3.5
$ws.Range("AX9").Value = 6.5
$ws.Range("AY9").Value = 15

# Row 11 (Racing Santander vs Burgos CF)
$ws.Range("G11").Value = 1.57
$ws.Range("H11").Value = 4.1
$ws.Range("I11").Value = 5.25
$ws.Range("M11").Value = 1.04
$ws.Range("N11").Value = 13
$ws.Range("Q11").Value = 1.7
$ws.Range("R11").Value = 2.1
$ws.Range("Z11").Value = 12
$ws.Range("AH11").Value = 15
$ws.Range("AO11").Value = 8
$ws.Range("AS11").Value = 101
$ws.Range("AX11").Value = 26
$ws.Range("AY11").Value = 29
$ws.Range("AZ11").Value = 81
